$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 3491.8823
$ws.Range("I62").Value = 3460.125
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 3460.125
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -2836.125
$ws.Range("N62").Value = -5248
# Row 65
$ws.Range("H65").Value = 3491.8823
$ws.Range("I65").Value = 3460.125
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 17300.625
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -14180.625
$ws.Range("N65").Value = -26240
# Row 70
$ws.Range("H70").Value = 4088.077
$ws.Range("J70").Value = 4303.909
$ws.Range("L70").Value = 12911.727
$ws.Range("N70").Value = -13451.727
# Row 73
$ws.Range("H73").Value = 4088.077
$ws.Range("J73").Value = 4303.909
$ws.Range("L73").Value = 12911.727
$ws.Range("N73").Value = -14783.727
# Row 98
$ws.Range("H98").Value = 1401
$ws.Range("I98").Value = 830.4286
$ws.Range("J98").Value = 1971.5714
$ws.Range("K98").Value = 830.4286
$ws.Range("L98").Value = 1971.5714
$ws.Range("M98").Value = 667.5714
$ws.Range("N98").Value = -4967.5714
# Row 100
$ws.Range("H100").Value = 2766.4443
$ws.Range("I100").Value = 2633
$ws.Range("K100").Value = 2633
$ws.Range("M100").Value = -2092
# Row 103
$ws.Range("H103").Value = 1192.5
$ws.Range("I103").Value = 1030
$ws.Range("J103").Value = 2005
$ws.Range("K103").Value = 3090
$ws.Range("L103").Value = 6015
$ws.Range("M103").Value = -2504
$ws.Range("N103").Value = -7187
# Row 122
$ws.Range("H122").Value = 1401
$ws.Range("I122").Value = 830.4286
$ws.Range("J122").Value = 1971.5714
$ws.Range("K122").Value = 2491.2858
$ws.Range("L122").Value = 5914.7142
$ws.Range("M122").Value = -41.28579999999965
$ws.Range("N122").Value = -10814.7142

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1988.9395
$ws.Range("I32").Value = 1847.9688
$ws.Range("K32").Value = 1847.9688
$ws.Range("M32").Value = -1560.9688
# Row 61
$ws.Range("H61").Value = 1681.1
$ws.Range("I61").Value = 1681.1
$ws.Range("K61").Value = 1681.1
$ws.Range("M61").Value = -1469.1
# Row 102
$ws.Range("H102").Value = 2130.3333
$ws.Range("I102").Value = 1496.75
$ws.Range("K102").Value = 1496.75
$ws.Range("M102").Value = 125.25
# Row 114
$ws.Range("H114").Value = 45000
$ws.Range("J114").Value = 45000
$ws.Range("L114").Value = 45000
$ws.Range("N114").Value = -53678
# Row 136
$ws.Range("H136").Value = 1681.1
$ws.Range("I136").Value = 1681.1
$ws.Range("K136").Value = 5043.299999999999
$ws.Range("M136").Value = -2493.299999999999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1217.3572
$ws.Range("I94").Value = 1249.8462
$ws.Range("K94").Value = 1249.8462
$ws.Range("M94").Value = -798.8462
# Row 105
$ws.Range("H105").Value = 2418.4614
$ws.Range("I105").Value = 2203.4167
$ws.Range("K105").Value = 2203.4167
$ws.Range("M105").Value = -456.4167000000002
# Row 107
$ws.Range("H107").Value = 1223.3846
$ws.Range("I107").Value = 1139.25
$ws.Range("K107").Value = 1139.25
$ws.Range("M107").Value = 780.75
# Row 134
$ws.Range("H134").Value = 8181.25
$ws.Range("I134").Value = 8181.25
$ws.Range("K134").Value = 24543.75
$ws.Range("M134").Value = -22008.75

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 116
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
# Row 122
$ws.Range("H122").Value = 2597.25
$ws.Range("I122").Value = 2597.25
$ws.Range("K122").Value = 7791.75
$ws.Range("M122").Value = -5341.75
# Row 132
$ws.Range("H132").Value = 2411.7334
$ws.Range("I132").Value = 2411.7334
$ws.Range("K132").Value = 7235.2002
$ws.Range("M132").Value = -4705.2002

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 128
$ws.Range("H128").Value = 340377.72
$ws.Range("I128").Value = 340377.72
$ws.Range("K128").Value = 1021133.16
$ws.Range("M128").Value = -1016153.16

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 96
$ws.Range("H96").Value = 50000
$ws.Range("J96").Value = 50000
$ws.Range("L96").Value = 50000
$ws.Range("N96").Value = -55492
# Row 97
$ws.Range("H97").Value = 707.61536
$ws.Range("I97").Value = 641.9
$ws.Range("J97").Value = 926.6667
$ws.Range("K97").Value = 641.9
$ws.Range("L97").Value = 926.6667
$ws.Range("M97").Value = -145.9
$ws.Range("N97").Value = -1918.6667
# Row 102
$ws.Range("H102").Value = 2034.4445
$ws.Range("I102").Value = 1404
$ws.Range("J102").Value = 4241
$ws.Range("K102").Value = 1404
$ws.Range("L102").Value = 4241
$ws.Range("M102").Value = 218
$ws.Range("N102").Value = -7485
# Row 113
$ws.Range("H113").Value = 4332.8887
$ws.Range("I113").Value = 4332.8887
$ws.Range("K113").Value = 4332.8887
$ws.Range("M113").Value = -2162.8887
# Row 122
$ws.Range("H122").Value = 2822.889
$ws.Range("I122").Value = 2599.389
$ws.Range("K122").Value = 7798.167
$ws.Range("M122").Value = -5348.167
# Row 132
$ws.Range("H132").Value = 2804.7144
$ws.Range("I132").Value = 2804.7144
$ws.Range("K132").Value = 8414.143199999999
$ws.Range("M132").Value = -5884.143199999999

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 1972.4445
$ws.Range("I40").Value = 1968
$ws.Range("K40").Value = 1968
$ws.Range("M40").Value = -1832
# Row 56
$ws.Range("H56").Value = 15000
$ws.Range("J56").Value = 15000
$ws.Range("L56").Value = 15000
$ws.Range("N56").Value = -16382
# Row 61
$ws.Range("H61").Value = 8701.6
$ws.Range("I61").Value = 7833.3335
$ws.Range("K61").Value = 7833.3335
$ws.Range("M61").Value = -7631.3335
# Row 113
$ws.Range("H113").Value = 8701.6
$ws.Range("I113").Value = 7833.3335
$ws.Range("K113").Value = 7833.3335
$ws.Range("M113").Value = -5663.3335
# Row 122
$ws.Range("H122").Value = 3789.2856
$ws.Range("I122").Value = 3504
$ws.Range("K122").Value = 10512
$ws.Range("M122").Value = -8062
# Row 132
$ws.Range("H132").Value = 13099.667
$ws.Range("I132").Value = 7579.4
$ws.Range("K132").Value = 22738.2
$ws.Range("M132").Value = -20208.2

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 3004.4
$ws.Range("I81").Value = 3227.111
$ws.Range("K81").Value = 6454.222
$ws.Range("M81").Value = -5393.222
# Row 84
$ws.Range("H84").Value = 3004.4
$ws.Range("I84").Value = 3227.111
$ws.Range("K84").Value = 32271.11
$ws.Range("M84").Value = -26967.11
